# Update cryptos list prices (col D) and volume(1h) percentages (col E).
# Source diff keeps plain-text cells (no explicit style); Excel auto-converts
# single-decimal numeric-looking strings to numbers on assignment, so those
# are written with a leading apostrophe (forces text) and then restyled back
# to "Normal" so no quote-prefix style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.068.50"
$ws.Range("E2").Value = "  +6.22%  "
$ws.Range("D3").Value = "2.331.57"
$ws.Range("E3").Value = "  +5.15%  "
$ws.Range("E4").Value = "  -0.88%  "
$ws.Range("D5").Value = "'305.44"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("D6").Value = "'97.86"
$ws.Range("E6").Value = "  +9.94%  "
$ws.Range("E7").Value = "  +5.01%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("E9").Value = "  +9.55%  "
$ws.Range("D10").Value = "'36.08"
$ws.Range("E10").Value = "  +7.83%  "
$ws.Range("D11").Value = "'0.0815"
$ws.Range("E11").Value = "  +4.88%  "
$ws.Range("E12").Value = "  +8.56%  "
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Value = "2.686.84"
$ws.Range("E14").Value = "  +4.99%  "
$ws.Range("D15").Value = "2.332.37"
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("E16").Value = "  +5.29%  "
$ws.Range("D17").Value = "'14.17"
$ws.Range("E17").Value = "  +8.15%  "
$ws.Range("D18").Value = "46.903.73"
$ws.Range("E18").Value = "  +6.43%  "
$ws.Range("E19").Value = "  +20.98%  "
$ws.Range("D20").Value = "0.0₃0958"
$ws.Range("E20").Value = "  +6.20%  "
$ws.Range("D21").Value = "'6.23"
$ws.Range("E21").Value = "  +4.24%  "
$ws.Range("D22").Value = "'68.15"
$ws.Range("E22").Value = "  +6.13%  "
$ws.Range("D23").Value = "'255.27"
$ws.Range("E23").Value = "  +9.76%  "
$ws.Range("D24").Value = "'2.99"
$ws.Range("E24").Value = "  +4.54%  "
$ws.Range("D25").Value = "'2.02"
$ws.Range("E25").Value = "  +5.63%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").Value = "'42.21"
$ws.Range("E27").Value = "  +17.11%  "
$ws.Range("D28").Value = "'2.32"
$ws.Range("E28").Value = "  +2.98%  "
$ws.Range("E29").Value = "  +6.29%  "
$ws.Range("D30").Value = "'20.39"
$ws.Range("E30").Value = "  +5.25%  "
$ws.Range("D31").Value = "'5.87"
$ws.Range("E31").Value = "  +5.12%  "
$ws.Range("D32").Value = "'0.0821"
$ws.Range("E32").Value = "  +9.19%  "
$ws.Range("D33").Value = "'147.64"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("E35").Value = "  +9.06%  "
$ws.Range("E36").Value = "  +6.75%  "
$ws.Range("E37").Value = "  +4.01%  "
$ws.Range("E38").Value = "  +5.49%  "
$ws.Range("E39").Value = "  +11.09%  "
$ws.Range("D40").Value = "'0.0311"
$ws.Range("E40").Value = "  +9.18%  "
$ws.Range("E41").Value = "  +7.09%  "
$ws.Range("D42").Value = "'14.11"
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").Value = "'1.98"
$ws.Range("E44").Value = "  +18.02%  "
$ws.Range("D45").Value = "'93.01"
$ws.Range("E45").Value = "  +18.53%  "
$ws.Range("D46").Value = "1.800.58"
$ws.Range("E46").Value = "  +3.65%  "
$ws.Range("D47").Value = "'74.81"
$ws.Range("E47").Value = "  +12.76%  "
$ws.Range("E48").Value = "  +8.88%  "
$ws.Range("E49").Value = "  +4.85%  "
$ws.Range("D50").Value = "'55.45"
$ws.Range("E50").Value = "  +7.04%  "
$ws.Range("D51").Value = "'4.86"
$ws.Range("E51").Value = "  +4.07%  "

# Reset quote-prefix style introduced above back to Normal so cells carry no explicit style,
# matching the unstyled D/E cells already in the sheet.
foreach ($addr in @("D5", "D6", "D8", "D10", "D11", "D17", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D40", "D42", "D43", "D44", "D45", "D47", "D50", "D51")) {
    $ws.Range($addr).Style = "Normal"
}
